$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.849.70'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '3.138.99'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'571.07"
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").Value = "'150.49"
$ws.Range("E6").Value = '  +3.44%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.135.75'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = '  +3.74%  '
$ws.Range("E10").Value = '  +4.77%  '
$ws.Range("D11").Value = "'6.18"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = "'0.502"
$ws.Range("E12").Value = '  +6.31%  '
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = '  +10.54%  '
$ws.Range("D14").Value = "'37.54"
$ws.Range("E14").Value = '  +6.49%  '
$ws.Range("D15").Value = '3.648.63'
$ws.Range("E15").Value = '  +1.86%  '
$ws.Range("D16").Value = '64.925.06'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = "'7.19"
$ws.Range("E17").Value = '  +6.03%  '
$ws.Range("D18").Value = '3.140.70'
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = "'510.95"
$ws.Range("E20").Value = '  +6.38%  '
$ws.Range("D21").Value = "'14.93"
$ws.Range("E21").Value = '  +7.07%  '
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = '  +8.17%  '
$ws.Range("D23").Value = "'15.46"
$ws.Range("E23").Value = '  +12.22%  '
$ws.Range("D24").Value = "'7.81"
$ws.Range("E24").Value = '  +3.13%  '
$ws.Range("D25").Value = "'85.05"
$ws.Range("E25").Value = '  +4.59%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  +3.38%  '
$ws.Range("D28").Value = "'8.71"
$ws.Range("E28").Value = '  +8.26%  '
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = '  +4.57%  '
$ws.Range("D30").Value = "'27.90"
$ws.Range("E30").Value = '  +6.24%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("D33").Value = "'2.64"
$ws.Range("E33").Value = '  +5.65%  '
$ws.Range("D34").Value = "'6.03"
$ws.Range("E34").Value = '  +7.71%  '
$ws.Range("D35").Value = "'6.58"
$ws.Range("E35").Value = '  +5.82%  '
$ws.Range("D36").Value = "'55.49"
$ws.Range("E36").Value = '  -0.51%  '
$ws.Range("D37").Value = "'477.06"
$ws.Range("E37").Value = '  +4.10%  '
$ws.Range("D38").Value = "'0.0423"
$ws.Range("E38").Value = '  +3.65%  '
$ws.Range("D39").Value = "'0.0858"
$ws.Range("E39").Value = '  +4.08%  '
$ws.Range("D40").Value = "'3.01"
$ws.Range("E40").Value = '  -1.79%  '
$ws.Range("D41").Value = '3.108.15'
$ws.Range("E41").Value = '  +4.47%  '
$ws.Range("D42").Value = "'8.60"
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("E43").Value = '  +3.84%  '
$ws.Range("D44").Value = "'0.291"
$ws.Range("E44").Value = '  +11.21%  '
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = '  +13.29%  '
$ws.Range("D46").Value = "'29.05"
$ws.Range("E46").Value = '  +3.60%  '
$ws.Range("D47").Value = '0.0₃0573'
$ws.Range("E47").Value = '  +10.71%  '
$ws.Range("E49").Value = '  +3.02%  '
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = '  +9.87%  '
$ws.Range("D51").Value = "'118.63"
$ws.Range("E51").Value = '  -1.91%  '
